$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-06-07 Saturday"

# Update each answer cell in the table, addressed by (row, col)
# so that duplicate text values (e.g. "11+19=30" appearing twice)
# are each replaced with their own distinct target value.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "71+10=81"
$t.Cell(1,2).Range.Text = "73-37=36"
$t.Cell(1,3).Range.Text = "2+36=38"
$t.Cell(1,4).Range.Text = "75-14=61"
$t.Cell(1,5).Range.Text = "22-2=20"
$t.Cell(2,1).Range.Text = "28-13=15"
$t.Cell(2,2).Range.Text = "34-4=30"
$t.Cell(2,3).Range.Text = "47-32=15"
$t.Cell(2,4).Range.Text = "99-69=30"
$t.Cell(2,5).Range.Text = "80+18=98"
$t.Cell(3,1).Range.Text = "47-31=16"
$t.Cell(3,2).Range.Text = "80+13=93"
$t.Cell(3,3).Range.Text = "74-14=60"
$t.Cell(3,4).Range.Text = "15-4=11"
$t.Cell(3,5).Range.Text = "63-1=62"
$t.Cell(4,1).Range.Text = "47+20=67"
$t.Cell(4,2).Range.Text = "99+0=99"
$t.Cell(4,3).Range.Text = "68+20=88"
$t.Cell(4,4).Range.Text = "51-25=26"
$t.Cell(4,5).Range.Text = "4+66=70"
$t.Cell(5,1).Range.Text = "21+33=54"
$t.Cell(5,2).Range.Text = "98-85=13"
$t.Cell(5,3).Range.Text = "12+51=63"
$t.Cell(5,4).Range.Text = "14+13=27"
$t.Cell(5,5).Range.Text = "65+2=67"
$t.Cell(6,1).Range.Text = "80-2=78"
$t.Cell(6,2).Range.Text = "59-0=59"
$t.Cell(6,3).Range.Text = "28+49=77"
$t.Cell(6,4).Range.Text = "48+28=76"
$t.Cell(6,5).Range.Text = "16+58=74"
$t.Cell(7,1).Range.Text = "13+67=80"
$t.Cell(7,2).Range.Text = "55-51=4"
$t.Cell(7,3).Range.Text = "97-50=47"
$t.Cell(7,4).Range.Text = "54+41=95"
$t.Cell(7,5).Range.Text = "43-5=38"
$t.Cell(8,1).Range.Text = "43+44=87"
$t.Cell(8,2).Range.Text = "56+42=98"
$t.Cell(8,3).Range.Text = "5+57=62"
$t.Cell(8,4).Range.Text = "57+40=97"
$t.Cell(8,5).Range.Text = "83-4=79"
$t.Cell(9,1).Range.Text = "97-41=56"
$t.Cell(9,2).Range.Text = "34+2=36"
$t.Cell(9,3).Range.Text = "33+12=45"
$t.Cell(9,4).Range.Text = "7+55=62"
$t.Cell(9,5).Range.Text = "4-1=3"
$t.Cell(10,1).Range.Text = "29+29=58"
$t.Cell(10,2).Range.Text = "0+99=99"
$t.Cell(10,3).Range.Text = "9+84=93"
$t.Cell(10,4).Range.Text = "45+36=81"
$t.Cell(10,5).Range.Text = "91-57=34"
$t.Cell(11,1).Range.Text = "83+6=89"
$t.Cell(11,2).Range.Text = "17+17=34"
$t.Cell(11,3).Range.Text = "47+46=93"
$t.Cell(11,4).Range.Text = "62-56=6"
$t.Cell(11,5).Range.Text = "58-29=29"
$t.Cell(12,1).Range.Text = "12+60=72"
$t.Cell(12,2).Range.Text = "94-9=85"
$t.Cell(12,3).Range.Text = "92-70=22"
$t.Cell(12,4).Range.Text = "85-80=5"
$t.Cell(12,5).Range.Text = "50+5=55"
$t.Cell(13,1).Range.Text = "27-4=23"
$t.Cell(13,2).Range.Text = "55+26=81"
$t.Cell(13,3).Range.Text = "83-83=0"
$t.Cell(13,4).Range.Text = "69-45=24"
$t.Cell(13,5).Range.Text = "75-65=10"
$t.Cell(14,1).Range.Text = "78-74=4"
$t.Cell(14,2).Range.Text = "44+18=62"
$t.Cell(14,3).Range.Text = "42+24=66"
$t.Cell(14,4).Range.Text = "48-21=27"
$t.Cell(14,5).Range.Text = "77-41=36"
$t.Cell(15,1).Range.Text = "33+25=58"
$t.Cell(15,2).Range.Text = "68+30=98"
$t.Cell(15,3).Range.Text = "52+40=92"
$t.Cell(15,4).Range.Text = "54+30=84"
$t.Cell(15,5).Range.Text = "29+47=76"
$t.Cell(16,1).Range.Text = "13+24=37"
$t.Cell(16,2).Range.Text = "78+4=82"
$t.Cell(16,3).Range.Text = "23+76=99"
$t.Cell(16,4).Range.Text = "93-87=6"
$t.Cell(16,5).Range.Text = "38+14=52"
$t.Cell(17,1).Range.Text = "42-4=38"
$t.Cell(17,2).Range.Text = "23+30=53"
$t.Cell(17,3).Range.Text = "87+7=94"
$t.Cell(17,4).Range.Text = "6+50=56"
$t.Cell(17,5).Range.Text = "10+10=20"
$t.Cell(18,1).Range.Text = "71+19=90"
$t.Cell(18,2).Range.Text = "13+77=90"
$t.Cell(18,3).Range.Text = "39+8=47"
$t.Cell(18,4).Range.Text = "36-11=25"
$t.Cell(18,5).Range.Text = "28-27=1"
$t.Cell(19,1).Range.Text = "85-30=55"
$t.Cell(19,2).Range.Text = "33+41=74"
$t.Cell(19,3).Range.Text = "54-50=4"
$t.Cell(19,4).Range.Text = "48-9=39"
$t.Cell(19,5).Range.Text = "74+24=98"
$t.Cell(20,1).Range.Text = "77-38=39"
$t.Cell(20,2).Range.Text = "64-63=1"
$t.Cell(20,3).Range.Text = "97-52=45"
$t.Cell(20,4).Range.Text = "88-11=77"
$t.Cell(20,5).Range.Text = "21+30=51"
